# Daily attendance processing - 2025-10-13 22:19:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C width: 23 -> 22 ---
$ws.Columns.Item(3).ColumnWidth = 21.17

# --- Recorded By email list reorderings ---
$ws.Range("G3").Value = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G4").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G12").Value = "wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G19").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G25").Value = "Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G26").Value = "rana.abozaid@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"
$ws.Range("G34").Value = "wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg, mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg"
$ws.Range("G41").Value = "Monica.Eshak@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, marina_atef@med.asu.edu.eg"
$ws.Range("G42").Value = "youstina.magdy@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, marina_atef@med.asu.edu.eg"

# --- Class Statistics block updates ---
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 3

# L9/L10 hold their percentages as literal text (not numbers), so force text
# entry with a leading apostrophe, then reapply the original (General, style 5)
# formatting from a same-styled neighbour to drop the resulting quote-prefix flag.
$ws.Range("L9").Value = "'25.0%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("L10").Value = "'36.6%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# --- Row 11: BIOCHEMISTRY LAB/CBL session for C1 now recorded ---
# Re-use the existing "Recorded" (green) formatting already applied to row 3
# by copying its format only (values are left untouched by PasteSpecial formats).
$ws.Range("A3:I3").Copy()
$ws.Range("A11:I11").PasteSpecial(-4122)
$ws.Range("C11").Value = "BIOCHEMISTRY LAB/CBL"
$ws.Range("G11").Value = "salma.elgendy.std@med.asu.edu.eg"
$ws.Range("H11").Value = "56/221"
$ws.Range("I11").Value = "Recorded"

# --- Row 33: BIOCHEMISTRY LAB/CBL session for C2 now recorded ---
$ws.Range("A25:I25").Copy()
$ws.Range("A33:I33").PasteSpecial(-4122)
$ws.Range("C33").Value = "BIOCHEMISTRY LAB/CBL"
$ws.Range("G33").Value = "salma.elgendy.std@med.asu.edu.eg"
$ws.Range("H33").Value = "35/246"
$ws.Range("I33").Value = "Recorded"

# --- Per-subject summary table (HISTOLOGY C1 row 15 / PHARMACOLOGY C1->C2 row 16) ---
$ws.Range("O15").Value = 5
$ws.Range("P15").Value = 1

$ws.Range("R15").Value = "'22.7%"
$ws.Range("N15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").Value = "'42.9%"
$ws.Range("N15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

$ws.Range("O16").Value = 6
$ws.Range("P16").Value = 2

$ws.Range("R16").Value = "'27.3%"
$ws.Range("N16").Copy()
$ws.Range("R16").PasteSpecial(-4122)

$ws.Range("S16").Value = "'31.3%"
$ws.Range("N16").Copy()
$ws.Range("S16").PasteSpecial(-4122)
